# Update "want to go" counts (column F) and one time-range text (column E)
# on the "展览" and "全部类型" worksheets, matching the refreshed scrape
# output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (row -> new F value) ----
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 9926
    3  = 418
    4  = 2518
    7  = 188
    9  = 730
    11 = 1221
    13 = 3070
    14 = 2319
    16 = 2013
    20 = 1572
    21 = 531
    23 = 211
    24 = 228
    25 = 41
    26 = 361
    28 = 343
    29 = 548
    30 = 43
    31 = 197
    32 = 1559
    33 = 263
    34 = 1597
    35 = 84
    36 = 389
    37 = 40
    38 = 420
    39 = 874
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}
$ws1.Range("E35").Value = "2024.05.02 10:00-05.02 22:00"

# ---- Sheet "全部类型" (row -> new F value) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 9926
    3  = 418
    4  = 2518
    9  = 188
    11 = 730
    13 = 1221
    15 = 3070
    16 = 2319
    18 = 2013
    22 = 1572
    23 = 531
    25 = 211
    26 = 228
    27 = 41
    28 = 361
    30 = 343
    31 = 548
    35 = 43
    36 = 197
    37 = 1559
    39 = 263
    40 = 1597
    41 = 84
    43 = 389
    44 = 40
    45 = 420
    46 = 874
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
$ws4.Range("E41").Value = "2024.05.02 10:00-05.02 22:00"
